# Apply updates to "团队_鑫e贷总授信完成情况" sheet data.
# The workbook's only sheet ("Sheet1") holds per-person counts in columns
# C (报表完成数), D (外拓双算数), F (调整数), H (本月 = D+E+F+G), and
# I (完成率 = H/B). This script updates the specific cells whose cached
# values changed, recomputing the dependent H and I values directly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4 (瞿逸程)
$ws.Range("C4").Value = 1
$ws.Range("D4").Value = 5
$ws.Range("H4").Value = 5

# Row 12 (曹倩云)
$ws.Range("C12").Value = 4
$ws.Range("D12").Value = 83
$ws.Range("F12").Value = 38
$ws.Range("H12").Value = 123
$ws.Range("I12").Value = 1.892307692307692

# Row 13 (黄旭)
$ws.Range("C13").Value = 10
$ws.Range("F13").Value = 36
$ws.Range("H13").Value = 145
$ws.Range("I13").Value = 2.230769230769231

# Row 16 (茅敏艳)
$ws.Range("C16").Value = 3
$ws.Range("F16").Value = 34
$ws.Range("H16").Value = 74
$ws.Range("I16").Value = 1.233333333333333

# Row 17 (钱潇伟)
$ws.Range("C17").Value = 2
$ws.Range("F17").Value = 26
$ws.Range("H17").Value = 61
$ws.Range("I17").Value = 1.016666666666667

# Row 18 (杨小东)
$ws.Range("C18").Value = 1
$ws.Range("F18").Value = 40
$ws.Range("H18").Value = 72
$ws.Range("I18").Value = 1.2
